$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 263
$ws.Range("C3").Value = 161470
$ws.Range("C4").Value = 152478
$ws.Range("C5").Value = 8992
$ws.Range("C8").Value = 64.42
